# Rename the sheet from "Sheet1" to "test".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "test"

# Turn A1 (which held "test1") into a mailto: hyperlink whose display text is
# the email address. This replaces the shared string "test1" with
# "supppyy@hotmail.com", pulls in the built-in "Hyperlink" cell style
# (underlined, theme color 10) for A1, and records the relationship used by
# the <hyperlinks> element.
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:supppyy@hotmail.com", "", "", "supppyy@hotmail.com") | Out-Null

# Widen column A so the new, longer text fits.
$ws.Columns("A").AutoFit() | Out-Null

# Leave the selection on B5, matching the saved view state.
$ws.Range("B5").Select() | Out-Null
